$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VehicleData")

# --- Update headers (row 1), columns W:AC ---
$ws.Range("W1").Value = "Mileage Rounded to Nearest 50,000"
$ws.Range("X1").Value = "Engine Size Rounded"
$ws.Range("Y1").Value = "Price Filter"
$ws.Range("Z1").Value = "Mileage Filter"
$ws.Range("AA1").Value = "Engine Size Filter"
$ws.Range("AB1").Value = "MPG Filter"
$ws.Range("AC1").Value = "Master Filter"

# --- Update formulas for data rows 2 and 3, columns W:AC ---
foreach ($r in 2..3) {
    $ws.Range("W$r").Formula = "=(ROUNDDOWN((AVERAGE(VehicleData!G$r))/(50000),0))*(50000)"
    $ws.Range("X$r").Formula = "=ROUND((AVERAGE(VehicleData!P$r))/(1000),1)"
    $ws.Range("Y$r").Formula = "=IF((AVERAGE(VehicleData!V$r))=(30000),0,1)"
    $ws.Range("Z$r").Formula = "=IF((AVERAGE(VehicleData!W$r))>(50000),0,1)"
    $ws.Range("AA$r").Formula = "=IF((AVERAGE(VehicleData!X$r))>(2.5),0,1)"
    $ws.Range("AB$r").Formula = "=IF((AVERAGE(VehicleData!Q$r))<(30),0,1)"
    $ws.Range("AC$r").Formula = "=IF((SUM(VehicleData!Y" + $r + ":AB" + $r + "))=(4),1,0)"
}

# --- Remove now-unused trailing columns AD:AF (Year Filter, MPG Filter(old), Master Filter(old)) ---
$ws.Columns("AD:AF").Delete()

# --- Remove the helper lookup sheets no longer referenced ---
$wb.Worksheets.Item("MileageBand").Delete()
$wb.Worksheets.Item("PriceBands").Delete()
